$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.252.45"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.861.32"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'0.7052"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'242.32"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.3136"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.07818"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").Value = "'24.28"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").Value = "'0.08003"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").Value = "1.859.08"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "'94.14"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'5.183"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'0.6977"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "'6.407"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "29.249.13"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'252.93"
$ws.Range("E19").Value = "  +4.39%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "2.108.87"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'7.554"
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("D24").Value = "'0.9989"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'0.1567"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").Value = "'9.001"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").Value = "'160.06"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "'4.309"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").Value = "'4.272"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'1.207"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'0.05284"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").Value = "'1.891"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "'0.7504"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "'2.707"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D39").Value = "1.246.13"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "'2.740"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "'111.26"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8983"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.132"
$ws.Range("E43").Value = "  -7.12%  "
$ws.Range("D44").Value = "'70.69"
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("D45").Value = "'0.9986"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "2.006.13"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "'1.786"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").Value = "'9.494"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "'0.4302"
$ws.Range("E51").Value = "  -2.06%  "
